$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep Text format so numeric-looking strings
# (e.g. "1.000", "28.402.99") are preserved exactly as text, matching
# the original inlineStr cell types.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.402.99'
$ws.Range("E2").Value = '  +3.97%  '
$ws.Range("D3").Value = '1.805.32'
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '316.29'
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5476'
$ws.Range("E7").Value = '  +4.42%  '
$ws.Range("D8").Value = '0.3847'
$ws.Range("E8").Value = '  +6.69%  '
$ws.Range("D9").Value = '0.07592'
$ws.Range("E9").Value = '  +2.84%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.127'
$ws.Range("E10").Value = '  +3.08%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = '42.25'
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").Value = '0.9995'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '21.18'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").Value = '6.190'
$ws.Range("E14").Value = '  +1.98%  '
$ws.Range("D15").Value = '7.391'
$ws.Range("E15").Value = '  +5.66%  '
$ws.Range("D16").Value = '1.806.37'
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("D17").Value = '92.08'
$ws.Range("E17").Value = '  +4.06%  '
$ws.Range("D18").Value = '0.00001071'
$ws.Range("E18").Value = '  +2.38%  '
$ws.Range("D19").Value = '0.06446'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '17.36'
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("D22").Value = '5.979'
$ws.Range("E22").Value = '  +2.21%  '
$ws.Range("D23").Value = '28.428.45'
$ws.Range("E23").Value = '  +3.73%  '
$ws.Range("D24").Value = '11.47'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("D26").Value = '159.17'
$ws.Range("E26").Value = '  +3.69%  '
$ws.Range("D27").Value = '20.68'
$ws.Range("E27").Value = '  +2.90%  '
$ws.Range("D28").Value = '2.400'
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("D29").Value = '2.014.74'
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("D30").Value = '123.96'
$ws.Range("E30").Value = '  +2.04%  '
$ws.Range("D31").Value = '1.125'
$ws.Range("E31").Value = '  +5.89%  '
$ws.Range("E32").Value = '  +4.03%  '
$ws.Range("D33").Value = '5.754'
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("D34").Value = '3.681'
$ws.Range("E34").Value = '  +2.16%  '
$ws.Range("D35").Value = '0.2313'
$ws.Range("E35").Value = '  +13.98%  '
$ws.Range("D36").Value = '0.06458'
$ws.Range("E36").Value = '  +7.80%  '
$ws.Range("E37").Value = '  +4.01%  '
$ws.Range("D38").Value = '8.842'
$ws.Range("E38").Value = '  +9.34%  '
$ws.Range("D39").Value = '5.163'
$ws.Range("E39").Value = '  +6.35%  '
$ws.Range("D40").Value = '11.64'
$ws.Range("E40").Value = '  +3.62%  '
$ws.Range("D41").Value = '0.6425'
$ws.Range("E41").Value = '  +4.33%  '
$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").Value = '0.9993'
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.161'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("D44").Value = '1.383'
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("D45").Value = '13.51'
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("D46").Value = '0.5989'
$ws.Range("E46").Value = '  +3.77%  '
$ws.Range("D47").Value = '3.684'
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("D48").Value = '126.91'
$ws.Range("E48").Value = '  +4.43%  '
$ws.Range("D49").Value = '1.988'
$ws.Range("E49").Value = '  +5.07%  '
$ws.Range("D50").Value = '1.148'
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").Value = '0.06902'
$ws.Range("E51").Value = '  +2.69%  '
